# Update crypto price/volume data per the Fri Aug 16 11:49:12 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.155.79'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '2.592.28'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.37'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.89'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").Value = '2.612.37'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.339'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '3.049.56'
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = '58.251.23'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.34'
$ws.Range("E16").Value = '  -3.51%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.568.97'
$ws.Range("E18").Value = '  -4.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.70'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.20'
$ws.Range("E21").Value = '  -2.42%  '
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.27'
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.403'
$ws.Range("E26").Value = '  -3.22%  '
$ws.Range("D27").Value = '2.718.62'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.04'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '0.0₃0739'
$ws.Range("E30").Value = '  -6.99%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("E32").Value = '  -6.11%  '
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.73'
$ws.Range("E34").Value = '  -0.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.38'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.99'
$ws.Range("E36").Value = '  -3.21%  '
$ws.Range("E37").Value = '  -4.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.853'
$ws.Range("E38").Value = '  -5.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.854'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.06'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.46'
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.52'
$ws.Range("E42").Value = '  -2.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.996'
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.605'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '270.58'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.68'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0955'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.70'
$ws.Range("E48").Value = '  -2.46%  '
$ws.Range("E49").Value = '  -2.40%  '
$ws.Range("D50").Value = '1.965.19'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.61'
$ws.Range("E51").Value = '  +0.61%  '
